# Aula 05 - Estrutura de Repeticao For
# Applies the text/formatting edits described by the target diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 6 - shape 2 ("Text Placeholder 2")
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$tr6 = $shp6.TextFrame.TextRange

# Paragraph 1: "A função range() retorna uma série de números consecutivos.
#               Por padrão, ela inicia no número 0 e é incrementada adicionando 1."
$para1 = $tr6.Paragraphs(1)

# 1a) "função range() " -> "função " + "range( ) " (bold/red run split)
$t1 = $para1.Text
$needle = "range() "
$rel = $t1.IndexOf($needle)
$abs = $para1.Start + $rel
$sub = $tr6.Characters($abs, $needle.Length)
$sub.Text = "range( ) "

# 1b) "...Por padrão, ela..." -> "...Por " + "padrão" (bold) + ", ela..."
$para1 = $tr6.Paragraphs(1)
$t1 = $para1.Text
$needle = "padrão"
$rel = $t1.IndexOf($needle)
$abs = $para1.Start + $rel
$sub = $tr6.Characters($abs, $needle.Length)
$sub.Font.Bold = $true

# Paragraph 3: "O comando range(4), por exemplo, ... A sintaxe da função range() é:"
$para3 = $tr6.Paragraphs(3)

# 1c) "...da função<nbsp>range() é:" -> "...da função<nbsp>" + "range" (bold) + "() é:"
$t3 = $para3.Text
$needle = "range"
$rel = $t3.IndexOf("range() é")
$abs = $para3.Start + $rel
$sub = $tr6.Characters($abs, $needle.Length)
$sub.Font.Bold = $true

# ---------------------------------------------------------------
# Slide 7 - shape 2 ("Text Placeholder 2")
# ---------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(2)
$tr7 = $shp7.TextFrame.TextRange

# "A função range() é utilizada..." -> "A função range( ) é utilizada..."
$para = $tr7.Paragraphs(5)
$t = $para.Text
$needle = "função range() "
$rel = $t.IndexOf($needle)
$abs = $para.Start + $rel
$sub = $tr7.Characters($abs, $needle.Length)
$sub.Text = "função range( ) "

# ---------------------------------------------------------------
# Slide 8 - shape 2: "Exemplo1 for/range:" -> "Exemplo 1 for/range:"
# (whole run replaced in one go, so no run split occurs)
# ---------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$tr8 = $shp8.TextFrame.TextRange
$para = $tr8.Paragraphs(1)
$needle = "Exemplo1 for/range:"
$rel = $para.Text.IndexOf($needle)
$abs = $para.Start + $rel
$sub = $tr8.Characters($abs, $needle.Length)
$sub.Text = "Exemplo 1 for/range:"

# ---------------------------------------------------------------
# Slide 9 - shape 2: "Exemplo2 for range:" -> "Exemplo 2 for range:"
#                     "Exemplo3 for range:" -> "Exemplo 3 for range:"
# (whole run replaced in one go, so no run split occurs)
# ---------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tr9 = $shp9.TextFrame.TextRange

$fullText = $tr9.Text
$needle = "Exemplo2 for range:"
$rel = $fullText.IndexOf($needle)
$abs = 1 + $rel
$sub = $tr9.Characters($abs, $needle.Length)
$sub.Text = "Exemplo 2 for range:"

$fullText = $tr9.Text
$needle = "Exemplo3 for range:"
$rel = $fullText.IndexOf($needle)
$abs = 1 + $rel
$sub = $tr9.Characters($abs, $needle.Length)
$sub.Text = "Exemplo 3 for range:"
